$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting rows 29:39 down to 30:40
$ws.Rows.Item(29).Insert()

# Fill in new row 29 data
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44522
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112027
$ws.Range("G29").Value = "Melón"
$ws.Range("H29").Value = "Tuna"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 23000
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = 23500
$ws.Range("N29").Value = "`$/caja 16 unidades"
$ws.Range("O29").Value = "Provincia de Copiapó"
$ws.Range("P29").Value = 1469
$ws.Range("Q29").Value = 16
$ws.Range("R29").Value = "Hortaliza"
